$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.379.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.86"
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6292"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07638"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "

# Row 9
$ws.Range("E9").Value = "  -0.43%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.11%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.837.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.66%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.00001095"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.57%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.005"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6785"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.73%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.090.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.67%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.125"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.414.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.440"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1390"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.371"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.17%  "

# Row 28
$ws.Range("E28").Value = "  +0.03%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.467"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.297"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.33%  "

# Row 31
$ws.Range("E31").Value = "  -1.07%  "

# Row 32
$ws.Range("E32").Value = "  -0.49%  "

# Row 33
$ws.Range("E33").Value = "  +0.64%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.850"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.156"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7090"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.68%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.585"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.772"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.228.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.92%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01798"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.441"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.74%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9072"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.53%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.000.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.25%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.145"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.53%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.041"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.683"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.38%  "
